# Scheduled-runner price/profit refresh for the Leves profit sheets.
# Updates currentAveragePrice / NQ / HQ, LevePrice NQ/HQ and LeveProfit
# NQ/HQ columns (H:N) for the rows whose market data changed.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 315.3
$ws.Range("I4").Value = 160.16667
$ws.Range("J4").Value = 548
$ws.Range("K4").Value = 160.16667
$ws.Range("L4").Value = 548
$ws.Range("M4").Value = -46.16667000000001
$ws.Range("N4").Value = -776

$ws.Range("H18").Value = 947.5
$ws.Range("I18").Value = 695
$ws.Range("J18").Value = 1200
$ws.Range("K18").Value = 695
$ws.Range("L18").Value = 1200
$ws.Range("M18").Value = -411
$ws.Range("N18").Value = -1768

$ws.Range("H75").Value = 24684.545
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 24684.545
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 24684.545
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -26556.545

$ws.Range("H78").Value = 24684.545
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 24684.545
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 74053.63499999999
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -83413.63499999999

$ws.Range("H95").Value = 47000
$ws.Range("J95").Value = 47000
$ws.Range("L95").Value = 47000
$ws.Range("N95").Value = -52492

$ws.Range("H100").Value = 2072.9167
$ws.Range("I100").Value = 1711.5
$ws.Range("K100").Value = 1711.5
$ws.Range("M100").Value = -1170.5

$ws.Range("H138").Value = 2314.6477
$ws.Range("I138").Value = 1409.0209
$ws.Range("J138").Value = 3401.4
$ws.Range("K138").Value = 4227.0627
$ws.Range("L138").Value = 10204.2
$ws.Range("M138").Value = 912.9372999999996
$ws.Range("N138").Value = -20484.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3763.3333
$ws.Range("I74").Value = 1038.1177
$ws.Range("J74").Value = 22294.8
$ws.Range("K74").Value = 1038.1177
$ws.Range("L74").Value = 22294.8
$ws.Range("M74").Value = -164.1177
$ws.Range("N74").Value = -24042.8

$ws.Range("H77").Value = 3763.3333
$ws.Range("I77").Value = 1038.1177
$ws.Range("J77").Value = 22294.8
$ws.Range("K77").Value = 5190.5885
$ws.Range("L77").Value = 111474
$ws.Range("M77").Value = -822.5884999999998
$ws.Range("N77").Value = -120210

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 21808.334
$ws.Range("J95").Value = 21808.334
$ws.Range("L95").Value = 21808.334
$ws.Range("N95").Value = -27300.334

$ws.Range("H134").Value = 1506.0571
$ws.Range("I134").Value = 987.7778
$ws.Range("K134").Value = 2963.3334
$ws.Range("M134").Value = -428.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 75.52941
$ws.Range("I7").Value = 50.5
$ws.Range("J7").Value = 97.77778000000001
$ws.Range("K7").Value = 50.5
$ws.Range("L7").Value = 97.77778000000001
$ws.Range("M7").Value = 62.5
$ws.Range("N7").Value = -323.77778

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 850.2162
$ws.Range("I4").Value = 214.28572
$ws.Range("J4").Value = 998.6
$ws.Range("K4").Value = 642.85716
$ws.Range("L4").Value = 2995.8
$ws.Range("M4").Value = -530.85716
$ws.Range("N4").Value = -3219.8

$ws.Range("H6").Value = 87.833336
$ws.Range("I6").Value = 95.40000000000001
$ws.Range("K6").Value = 286.2
$ws.Range("M6").Value = -173.2

$ws.Range("H40").Value = 449.3846
$ws.Range("I40").Value = 148.4
$ws.Range("J40").Value = 637.5
$ws.Range("K40").Value = 593.6
$ws.Range("L40").Value = 2550
$ws.Range("M40").Value = -524.6
$ws.Range("N40").Value = -2688

$ws.Range("H131").Value = 948.88
$ws.Range("I131").Value = 700
$ws.Range("J131").Value = 996.2857
$ws.Range("K131").Value = 2100
$ws.Range("L131").Value = 2988.8571
$ws.Range("M131").Value = 2940
$ws.Range("N131").Value = -13068.8571

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 40.157894
$ws.Range("I2").Value = 31.222221
$ws.Range("J2").Value = 48.2
$ws.Range("K2").Value = 31.222221
$ws.Range("L2").Value = 48.2
$ws.Range("M2").Value = 81.777779
$ws.Range("N2").Value = -274.2

$ws.Range("H80").Value = 2943.1
$ws.Range("I80").Value = 2604.1667
$ws.Range("J80").Value = 3451.5
$ws.Range("K80").Value = 2604.1667
$ws.Range("L80").Value = 3451.5
$ws.Range("M80").Value = -1606.1667
$ws.Range("N80").Value = -5447.5

$ws.Range("H83").Value = 2943.1
$ws.Range("I83").Value = 2604.1667
$ws.Range("J83").Value = 3451.5
$ws.Range("K83").Value = 13020.8335
$ws.Range("L83").Value = 17257.5
$ws.Range("M83").Value = -8028.833500000001
$ws.Range("N83").Value = -27241.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 431
$ws.Range("I22").Value = 240
$ws.Range("J22").Value = 478.75
$ws.Range("K22").Value = 240
$ws.Range("L22").Value = 478.75
$ws.Range("M22").Value = 55
$ws.Range("N22").Value = -1068.75

$ws.Range("H27").Value = 431
$ws.Range("I27").Value = 240
$ws.Range("J27").Value = 478.75
$ws.Range("K27").Value = 240
$ws.Range("L27").Value = 478.75
$ws.Range("M27").Value = -133
$ws.Range("N27").Value = -692.75

$ws.Range("H82").Value = 1360.3636
$ws.Range("I82").Value = 1110.7
$ws.Range("J82").Value = 1744.4615
$ws.Range("K82").Value = 1110.7
$ws.Range("L82").Value = 1744.4615
$ws.Range("M82").Value = -749.7
$ws.Range("N82").Value = -2466.4615

$ws.Range("H85").Value = 1360.3636
$ws.Range("I85").Value = 1110.7
$ws.Range("J85").Value = 1744.4615
$ws.Range("K85").Value = 1110.7
$ws.Range("L85").Value = 1744.4615
$ws.Range("M85").Value = 137.3
$ws.Range("N85").Value = -4240.461499999999

$ws.Range("H93").Value = 1005.6071
$ws.Range("J93").Value = 1874.0769
$ws.Range("L93").Value = 1874.0769
$ws.Range("N93").Value = -4370.0769

$ws.Range("H132").Value = 15597.594
$ws.Range("I132").Value = 8439.789000000001
$ws.Range("J132").Value = 26059
$ws.Range("K132").Value = 25319.367
$ws.Range("L132").Value = 78177
$ws.Range("M132").Value = -22789.367
$ws.Range("N132").Value = -83237

$ws.Range("H136").Value = 5431.6772
$ws.Range("I136").Value = 1611.4445
$ws.Range("J136").Value = 10721.23
$ws.Range("K136").Value = 4834.333500000001
$ws.Range("L136").Value = 32163.69
$ws.Range("M136").Value = -2284.333500000001
$ws.Range("N136").Value = -37263.69

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 33600
$ws.Range("J109").Value = 33600
$ws.Range("L109").Value = 33600
$ws.Range("N109").Value = -36374

$ws.Range("H133").Value = 32315
$ws.Range("J133").Value = 32315
$ws.Range("L133").Value = 32315
$ws.Range("N133").Value = -42435
